# Update the dSF (column F) values for a set of rows, per the re-pulled /
# recalculated data described in the commit message "repull data, push
# all data, mean calculation".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = -6
    9  = 10
    13 = 1
    15 = 0
    19 = -2
    20 = 1
    21 = -4
    24 = 1
    26 = 1
    29 = -2
    36 = -1
    42 = 3
    49 = 2
    55 = 0
    60 = -1
    61 = 8
    65 = 3
    72 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
